$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells keep their exact text representation
# (e.g. trailing zeros, thousand-separator dots) instead of being
# auto-converted to numbers by Excels value parser.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '65.542.35'
$ws.Range("E2").Value = '  -0.54%  '

# Row 3
$ws.Range("D3").Value = '3.276.19'
$ws.Range("E3").Value = '  -0.95%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '575.08'
$ws.Range("E5").Value = '  +3.20%  '

# Row 6
$ws.Range("D6").Value = '182.08'
$ws.Range("E6").Value = '  -2.33%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = '3.270.06'
$ws.Range("E8").Value = '  -0.87%  '

# Row 9
$ws.Range("D9").Value = '0.568'
$ws.Range("E9").Value = '  -2.75%  '

# Row 10
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  -5.66%  '

# Row 11
$ws.Range("D11").Value = '0.568'
$ws.Range("E11").Value = '  -2.52%  '

# Row 12
$ws.Range("D12").Value = '45.98'
$ws.Range("E12").Value = '  -3.26%  '

# Row 13
$ws.Range("D13").Value = '0.0000262'
$ws.Range("E13").Value = '  -2.82%  '

# Row 14
$ws.Range("D14").Value = '3.797.18'
$ws.Range("E14").Value = '  -0.99%  '

# Row 15
$ws.Range("D15").Value = '8.35'
$ws.Range("E15").Value = '  -2.98%  '

# Row 16
$ws.Range("D16").Value = '612.79'
$ws.Range("E16").Value = '  -2.53%  '

# Row 17
$ws.Range("D17").Value = '65.635.81'
$ws.Range("E17").Value = '  -0.44%  '

# Row 18
$ws.Range("E18").Value = '  +0.29%  '

# Row 19
$ws.Range("D19").Value = '17.67'
$ws.Range("E19").Value = '  -2.40%  '

# Row 20
$ws.Range("D20").Value = '3.280.39'
$ws.Range("E20").Value = '  -0.92%  '

# Row 21
$ws.Range("D21").Value = '10.84'
$ws.Range("E21").Value = '  -3.23%  '

# Row 22
$ws.Range("D22").Value = '0.882'
$ws.Range("E22").Value = '  -2.57%  '

# Row 23
$ws.Range("D23").Value = '17.88'
$ws.Range("E23").Value = '  -1.45%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '4.92'
$ws.Range("E24").Value = '  -0.84%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '98.15'
$ws.Range("E25").Value = '  -3.13%  '

# Row 26
$ws.Range("D26").Value = '3.93'
$ws.Range("E26").Value = '  -0.13%  '

# Row 27
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").Value = '  -1.24%  '

# Row 28
$ws.Range("D28").Value = '9.40'
$ws.Range("E28").Value = '  -1.26%  '

# Row 29
$ws.Range("D29").Value = '30.67'
$ws.Range("E29").Value = '  +1.55%  '

# Row 30
$ws.Range("D30").Value = '8.32'
$ws.Range("E30").Value = '  -3.70%  '

# Row 31
$ws.Range("D31").Value = '6.44'
$ws.Range("E31").Value = '  +1.33%  '

# Row 32
$ws.Range("D32").Value = '3.72'
$ws.Range("E32").Value = '  -6.82%  '

# Row 33
$ws.Range("D33").Value = '10.79'
$ws.Range("E33").Value = '  -2.51%  '

# Row 34
$ws.Range("D34").Value = '541.97'
$ws.Range("E34").Value = '  -2.36%  '

# Row 35
$ws.Range("D35").Value = '3.784.26'
$ws.Range("E35").Value = '  -1.57%  '

# Row 36
$ws.Range("D36").Value = '0.102'
$ws.Range("E36").Value = '  -2.41%  '

# Row 37
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("D38").Value = '55.92'
$ws.Range("E38").Value = '  -2.72%  '

# Row 39
$ws.Range("D39").Value = '0.127'
$ws.Range("E39").Value = '  -0.99%  '

# Row 40
$ws.Range("D40").Value = '32.34'
$ws.Range("E40").Value = '  -3.98%  '

# Row 41
$ws.Range("E41").Value = '  +3.43%  '

# Row 42
$ws.Range("D42").Value = '3.11'
$ws.Range("E42").Value = '  -4.60%  '

# Row 43
$ws.Range("D43").Value = '0.0₃0674'
$ws.Range("E43").Value = '  -7.91%  '

# Row 44
$ws.Range("D44").Value = '2.56'
$ws.Range("E44").Value = '  -4.55%  '

# Row 45
$ws.Range("D45").Value = '0.328'
$ws.Range("E45").Value = '  -1.52%  '

# Row 46
$ws.Range("D46").Value = '0.0402'
$ws.Range("E46").Value = '  -3.51%  '

# Row 47
$ws.Range("E47").Value = '  -8.23%  '

# Row 48
$ws.Range("E48").Value = '  +0.36%  '

# Row 49
$ws.Range("D49").Value = '0.125'
$ws.Range("E49").Value = '  -2.56%  '

# Row 50
$ws.Range("D50").Value = '2.48'
$ws.Range("E50").Value = '  -4.59%  '

# Row 51
$ws.Range("D51").Value = '127.86'
$ws.Range("E51").Value = '  +4.64%  '
